$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.285.33"
$ws.Range("E2").Value = "  -2.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.370.46"
$ws.Range("E3").Value = "  -2.51%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "500.21"
$ws.Range("E5").Value = "  -1.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.73"
$ws.Range("E6").Value = "  -2.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.545"
$ws.Range("E8").Value = "  -3.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.374.65"
$ws.Range("E9").Value = "  -4.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0971"
$ws.Range("E10").Value = "  -1.70%  "
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.61"
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.790.21"
$ws.Range("E14").Value = "  -2.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.199.26"
$ws.Range("E15").Value = "  -2.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.50"
$ws.Range("E16").Value = "  -2.97%  "
$ws.Range("E17").Value = "  -2.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.351.95"
$ws.Range("E18").Value = "  -3.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.05"
$ws.Range("E19").Value = "  -3.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.01"
$ws.Range("E20").Value = "  -3.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "307.32"
$ws.Range("E21").Value = "  -3.21%  "
$ws.Range("E22").Value = "  -3.75%  "
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.25"
$ws.Range("E24").Value = "  -0.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.996"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.370"
$ws.Range("E26").Value = "  -4.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.148"
$ws.Range("E27").Value = "  -4.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.28"
$ws.Range("E28").Value = "  -5.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.32"
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0714"
$ws.Range("E30").Value = "  -4.34%  "
$ws.Range("E31").Value = "  -4.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.76"
$ws.Range("E33").Value = "  -8.20%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.08"
$ws.Range("E34").Value = "  -7.09%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.56"
$ws.Range("E36").Value = "  -3.45%  "
$ws.Range("E37").Value = "  -7.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.78"
$ws.Range("E38").Value = "  -3.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.06"
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.795"
$ws.Range("E40").Value = "  -3.78%  "
$ws.Range("E41").Value = "  -4.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "130.39"
$ws.Range("E42").Value = "  -4.32%  "
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.77"
$ws.Range("E44").Value = "  -6.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.565"
$ws.Range("E45").Value = "  -2.45%  "
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "241.18"
$ws.Range("E47").Value = "  -8.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0482"
$ws.Range("E48").Value = "  -4.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0208"
$ws.Range("E49").Value = "  -4.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.11"
$ws.Range("E50").Value = "  -2.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.55"
$ws.Range("E51").Value = "  -3.86%  "
